# Auto-generated PowerShell COM-interop edit script.
# Commit: "Updated cryptos list on Tue Jan 30 15:06:56 UTC 2024 with GitHub Actions"
#
# The sheet stores every data cell as a literal (inline) string - including
# "Price" values that look numeric (e.g. "36.00", "43.424.12"). Plain
# `Range.Value = "..."` assignment lets Excel's COM layer auto-coerce
# plain-looking numeric text into a real Number (dropping trailing zeros /
# normalizing), which would corrupt cells such as D10 ("36.00" -> 36) or
# silently retype others. To keep every written cell a genuine text value
# (matching the original inlineStr cells), we force text NumberFormat,
# assign the literal string, then restore the cell style to "Normal" so no
# stray formatting is introduced.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}


Set-TextValue "D2" "43.424.12"
$ws.Range("E2").Value = "  +3.64%  "
Set-TextValue "D3" "2.310.29"
$ws.Range("E3").Value = "  +3.08%  "
$ws.Range("B5").Value = "BNB"
$ws.Range("C5").Value = "https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb"
Set-TextValue "D5" "308.49"
$ws.Range("E5").Value = "  +0.81%  "
$ws.Range("B6").Value = "Solana"
$ws.Range("C6").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
Set-TextValue "D6" "105.24"
$ws.Range("E6").Value = "  +9.80%  "
$ws.Range("E7").Value = "  +0.62%  "
$ws.Range("E8").Value = "  +0.01%  "
Set-TextValue "D9" "0.518"
$ws.Range("E9").Value = "  +6.21%  "
Set-TextValue "D10" "36.00"
$ws.Range("E10").Value = "  +4.39%  "
Set-TextValue "D11" "52.84"
$ws.Range("E11").Value = "  +2.78%  "
Set-TextValue "D12" "0.0811"
$ws.Range("E12").Value = "  -0.29%  "
$ws.Range("E13").Value = "  -1.05%  "
Set-TextValue "D14" "6.97"
$ws.Range("E14").Value = "  +3.59%  "
Set-TextValue "D15" "2.668.80"
$ws.Range("E15").Value = "  +3.09%  "
Set-TextValue "D16" "15.11"
$ws.Range("E16").Value = "  +5.47%  "
Set-TextValue "D17" "2.306.73"
$ws.Range("E17").Value = "  +2.38%  "
$ws.Range("E18").Value = "  +2.90%  "
Set-TextValue "D19" "43.364.93"
$ws.Range("E19").Value = "  +3.73%  "
Set-TextValue "D20" "0.0₃0922"
$ws.Range("E20").Value = "  +2.54%  "
Set-TextValue "D21" "11.86"
$ws.Range("E21").Value = "  -1.90%  "
$ws.Range("E22").Value = "  +4.99%  "
Set-TextValue "D23" "67.99"
$ws.Range("E23").Value = "  +1.45%  "
Set-TextValue "D24" "240.44"
$ws.Range("E24").Value = "  +2.38%  "
Set-TextValue "D25" "2.03"
$ws.Range("E25").Value = "  +5.21%  "
Set-TextValue "D26" "2.60"
$ws.Range("E26").Value = "  +1.73%  "
$ws.Range("E27").Value = "  -0.10%  "
Set-TextValue "D28" "24.99"
$ws.Range("E28").Value = "  +8.02%  "
Set-TextValue "D29" "2.21"
$ws.Range("E29").Value = "  +4.98%  "
Set-TextValue "D30" "36.28"
$ws.Range("E30").Value = "  -3.19%  "
Set-TextValue "D31" "9.58"
$ws.Range("E31").Value = "  +1.73%  "
Set-TextValue "D32" "162.46"
$ws.Range("E32").Value = "  -2.51%  "
$ws.Range("E33").Value = "  +1.78%  "
$ws.Range("E34").Value = "  -0.15%  "
Set-TextValue "D35" "18.29"
$ws.Range("E35").Value = "  +5.67%  "
$ws.Range("E36").Value = "  +6.69%  "
$ws.Range("E37").Value = "  +2.25%  "
Set-TextValue "D38" "4.60"
$ws.Range("E38").Value = "  +14.28%  "
Set-TextValue "D39" "3.01"
$ws.Range("B40").Value = "ARBITRUM"
$ws.Range("C40").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
Set-TextValue "D40" "1.86"
$ws.Range("E40").Value = "  +5.04%  "
$ws.Range("B41").Value = "Kaspa"
$ws.Range("C41").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
Set-TextValue "D41" "0.105"
$ws.Range("E41").Value = "  +3.15%  "
$ws.Range("E42").Value = "  +0.50%  "
Set-TextValue "D43" "2.48"
$ws.Range("E43").Value = "  +14.49%  "
$ws.Range("E44").Value = "  +3.41%  "
Set-TextValue "D45" "1.963.59"
$ws.Range("E45").Value = "  +1.49%  "
Set-TextValue "D46" "18.91"
$ws.Range("E46").Value = "  +2.59%  "
Set-TextValue "D47" "3.06"
$ws.Range("E47").Value = "  +6.62%  "
Set-TextValue "D48" "10.23"
$ws.Range("E48").Value = "  +6.61%  "
Set-TextValue "D49" "57.87"
$ws.Range("E49").Value = "  +7.75%  "
$ws.Range("E50").Value = "  +2.65%  "
$ws.Range("E51").Value = "  +8.10%  "
